{"js": "// Replace the worksheet date and each \"a\u00f7b=\" problem with its new value.\n// The replacements are positional (document order) because some source\n// values (e.g. \"49\u00f73=\") repeat with different targets, so a naive global\n// find/replace would corrupt the mapping.\nconst replacements = [\n  \"2023-09-23 Saturday\",\n  \"71\u00f76=\",\n  \"81\u00f77=\",\n  \"45\u00f75=\",\n  \"80\u00f75=\",\n  \"27\u00f73=\",\n  \"69\u00f74=\",\n  \"40\u00f77=\",\n  \"64\u00f72=\",\n  \"78\u00f72=\",\n  \"79\u00f75=\",\n  \"71\u00f74=\",\n  \"20\u00f73=\",\n  \"31\u00f72=\",\n  \"39\u00f72=\",\n  \"97\u00f77=\",\n  \"48\u00f73=\",\n  \"61\u00f78=\",\n  \"52\u00f75=\",\n  \"70\u00f75=\",\n  \"33\u00f79=\",\n  \"52\u00f75=\",\n  \"44\u00f73=\",\n  \"70\u00f78=\",\n  \"29\u00f79=\",\n  \"12\u00f72=\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet replaceIdx = 0;\nfor (let i = 0; i < paragraphs.items.length && replaceIdx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === null || text === undefined || text === \"\") {\n    continue;\n  }\n  para.insertText(replacements[replaceIdx], \"Replace\");\n  replaceIdx++;\n}\n\nawait context.sync();\n", "ps1": "# Replace the worksheet date and each \"a\u00f7b=\" problem with its new value.\n# The replacements are positional (document order) because some source\n# values (e.g. \"49\u00f73=\") repeat with different targets, so a naive global\n# find/replace would corrupt the mapping. We walk the document once,\n# searching forward from the end of the previous hit for the Nth old\n# value and swapping in the Nth new value.\n\n$pairs = @(\n    @{ Old = \"2023-09-22 Friday\"; New = \"2023-09-23 Saturday\" },\n    @{ Old = \"40\u00f73=\"; New = \"71\u00f76=\" },\n    @{ Old = \"80\u00f79=\"; New = \"81\u00f77=\" },\n    @{ Old = \"30\u00f72=\"; New = \"45\u00f75=\" },\n    @{ Old = \"66\u00f73=\"; New = \"80\u00f75=\" },\n    @{ Old = \"69\u00f77=\"; New = \"27\u00f73=\" },\n    @{ Old = \"23\u00f75=\"; New = \"69\u00f74=\" },\n    @{ Old = \"73\u00f74=\"; New = \"40\u00f77=\" },\n    @{ Old = \"90\u00f76=\"; New = \"64\u00f72=\" },\n    @{ Old = \"61\u00f78=\"; New = \"78\u00f72=\" },\n    @{ Old = \"80\u00f74=\"; New = \"79\u00f75=\" },\n    @{ Old = \"35\u00f79=\"; New = \"71\u00f74=\" },\n    @{ Old = \"96\u00f78=\"; New = \"20\u00f73=\" },\n    @{ Old = \"36\u00f73=\"; New = \"31\u00f72=\" },\n    @{ Old = \"24\u00f73=\"; New = \"39\u00f72=\" },\n    @{ Old = \"70\u00f74=\"; New = \"97\u00f77=\" },\n    @{ Old = \"42\u00f75=\"; New = \"48\u00f73=\" },\n    @{ Old = \"46\u00f78=\"; New = \"61\u00f78=\" },\n    @{ Old = \"49\u00f73=\"; New = \"52\u00f75=\" },\n    @{ Old = \"38\u00f76=\"; New = \"70\u00f75=\" },\n    @{ Old = \"29\u00f73=\"; New = \"33\u00f79=\" },\n    @{ Old = \"58\u00f76=\"; New = \"52\u00f75=\" },\n    @{ Old = \"97\u00f73=\"; New = \"44\u00f73=\" },\n    @{ Old = \"49\u00f73=\"; New = \"70\u00f78=\" },\n    @{ Old = \"96\u00f76=\"; New = \"29\u00f79=\" },\n    @{ Old = \"52\u00f76=\"; New = \"12\u00f72=\" }\n)\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n\nforeach ($pair in $pairs) {\n    $rng.Find.ClearFormatting()\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWildcards = $false\n    $found = $rng.Find.Execute($pair.Old)\n    if ($found) {\n        $rng.Text = $pair.New\n        $rng.Collapse(0)\n    }\n}\n"}
